$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new columns I (I0) and J (IF) ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header formatting (bold font, thin border, center/top align)
# used by the existing header cells (e.g. H1) by copying its format.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Data rows 2..62: values for columns I and J ---
$iValues = @(5,6,6,7,6,7,9,8,7,8,7,8,8,8,8,9,8,7,7,8,8,8,8,7,8,7,8,13,5,7,8,7,7,5,6,8,9,8,8,7,6,7,7,7,7,6,6,6,8,7,7,7,7,8,8,7,8,7,8,7,5)
$jValues = @(7,8,7,7,8,8,9,8,8,9,8,8,8,8,8,9,8,8,8,9,9,8,9,7,8,8,8,13,5,7,8,7,7,5,6,8,9,8,8,7,6,7,7,8,8,6,6,7,8,7,8,7,8,8,8,8,8,8,8,7,5)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
